$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 16, shifting the "quake" row and everything below it down by one.
$ws.Rows("16:16").Insert()

# Populate the new row 16 with the "stinger" spell data.
$ws.Range("A16").Value = "stinger"
$ws.Range("B16").Value = "attack"
$ws.Range("C16").Value = "High Dmg Critical Hit"
$ws.Range("D16").Value = "anim_spell_stab"
$ws.Range("E16").Value = "sfx_stab"
$ws.Range("F16").Value = 100
$ws.Range("G16").Value = 1
$ws.Range("H16").Value = 9999
$ws.Range("I16").Value = 9999
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = "MT"
$ws.Range("L16").Value = "ADSHE"
$ws.Range("M16").Value = 100

# Match the saved view state: M16 selected (and the window scrolled so row 7
# is at the top, to the extent the host lets us drive that).
$ws.Range("M16").Select()
try { $excel.ActiveWindow.ScrollRow = 7 } catch {}
